$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 383772.53
$ws.Range("I2").Value = 808323.3
$ws.Range("J2").Value = 1676.8
$ws.Range("K2").Value = 808323.3
$ws.Range("L2").Value = 1676.8
$ws.Range("M2").Value = -808210.3
$ws.Range("N2").Value = -1902.8
$ws.Range("H40").Value = 4723.6665
$ws.Range("J40").Value = 6131.3335
$ws.Range("L40").Value = 6131.3335
$ws.Range("N40").Value = -6481.3335
$ws.Range("H41").Value = 786.625
$ws.Range("I41").Value = 1059.4
$ws.Range("K41").Value = 1059.4
$ws.Range("M41").Value = -619.4000000000001
$ws.Range("H53").Value = 915.0769
$ws.Range("I53").Value = 723.1111
$ws.Range("J53").Value = 1347
$ws.Range("K53").Value = 723.1111
$ws.Range("L53").Value = 1347
$ws.Range("M53").Value = -86.11109999999996
$ws.Range("N53").Value = -2621
$ws.Range("H76").Value = 5081.6665
$ws.Range("I76").Value = 4998.125
$ws.Range("J76").Value = 5750
$ws.Range("K76").Value = 4998.125
$ws.Range("L76").Value = 5750
$ws.Range("M76").Value = -4683.125
$ws.Range("N76").Value = -6380
$ws.Range("H79").Value = 5081.6665
$ws.Range("I79").Value = 4998.125
$ws.Range("J79").Value = 5750
$ws.Range("K79").Value = 4998.125
$ws.Range("L79").Value = 5750
$ws.Range("M79").Value = -3906.125
$ws.Range("N79").Value = -7934
$ws.Range("H86").Value = 69589
$ws.Range("I86").Value = 86536.25
$ws.Range("J86").Value = 1800
$ws.Range("K86").Value = 86536.25
$ws.Range("L86").Value = 1800
$ws.Range("M86").Value = -85413.25
$ws.Range("N86").Value = -4046
$ws.Range("H89").Value = 69589
$ws.Range("I89").Value = 86536.25
$ws.Range("J89").Value = 1800
$ws.Range("K89").Value = 432681.25
$ws.Range("L89").Value = 9000
$ws.Range("M89").Value = -427065.25
$ws.Range("N89").Value = -20232
$ws.Range("H98").Value = 3390.2856
$ws.Range("I98").Value = 2460.5
$ws.Range("J98").Value = 4630
$ws.Range("K98").Value = 2460.5
$ws.Range("L98").Value = 4630
$ws.Range("M98").Value = -962.5
$ws.Range("N98").Value = -7626
$ws.Range("H106").Value = 3551.5454
$ws.Range("I106").Value = 3551.5454
$ws.Range("K106").Value = 3551.5454
$ws.Range("M106").Value = -2920.5454
$ws.Range("H122").Value = 3390.2856
$ws.Range("I122").Value = 2460.5
$ws.Range("J122").Value = 4630
$ws.Range("K122").Value = 7381.5
$ws.Range("L122").Value = 13890
$ws.Range("M122").Value = -4931.5
$ws.Range("N122").Value = -18790
$ws.Range("H129").Value = 8078.3335
$ws.Range("I129").Value = 2674.3333
$ws.Range("K129").Value = 8022.999899999999
$ws.Range("M129").Value = -3022.999899999999

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 3353.8
$ws.Range("I63").Value = 2949.4285
$ws.Range("J63").Value = 4297.3335
$ws.Range("K63").Value = 2949.4285
$ws.Range("L63").Value = 4297.3335
$ws.Range("M63").Value = -2263.4285
$ws.Range("N63").Value = -5669.3335
$ws.Range("H66").Value = 3353.8
$ws.Range("I66").Value = 2949.4285
$ws.Range("J66").Value = 4297.3335
$ws.Range("K66").Value = 14747.1425
$ws.Range("L66").Value = 21486.6675
$ws.Range("M66").Value = -11315.1425
$ws.Range("N66").Value = -28350.6675
$ws.Range("H88").Value = 2264.6875
$ws.Range("I88").Value = 2270
$ws.Range("J88").Value = 2261.5
$ws.Range("K88").Value = 2270
$ws.Range("L88").Value = 2261.5
$ws.Range("M88").Value = -1864
$ws.Range("N88").Value = -3073.5
$ws.Range("H91").Value = 2264.6875
$ws.Range("I91").Value = 2270
$ws.Range("J91").Value = 2261.5
$ws.Range("K91").Value = 2270
$ws.Range("L91").Value = 2261.5
$ws.Range("M91").Value = -866
$ws.Range("N91").Value = -5069.5

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H35").Value = 0
$ws.Range("I35").Value = 0
$ws.Range("K35").Value = 0
$ws.Range("H86").Value = 2762.375
$ws.Range("J86").Value = 3066.5
$ws.Range("L86").Value = 3066.5
$ws.Range("N86").Value = -5312.5
$ws.Range("H89").Value = 2762.375
$ws.Range("J89").Value = 3066.5
$ws.Range("L89").Value = 15332.5
$ws.Range("N89").Value = -26564.5
$ws.Range("H105").Value = 4541.1816
$ws.Range("I105").Value = 4425.3125
$ws.Range("J105").Value = 4850.1665
$ws.Range("K105").Value = 4425.3125
$ws.Range("L105").Value = 4850.1665
$ws.Range("M105").Value = -2678.3125
$ws.Range("N105").Value = -8344.166499999999
$ws.Range("H134").Value = 2154.5918
$ws.Range("I134").Value = 1866.2307
$ws.Range("J134").Value = 3279.2
$ws.Range("K134").Value = 5598.6921
$ws.Range("L134").Value = 9837.599999999999
$ws.Range("M134").Value = -3063.6921
$ws.Range("N134").Value = -14907.6

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 3800
$ws.Range("I62").Value = 3733.3333
$ws.Range("J62").Value = 4000
$ws.Range("K62").Value = 3733.3333
$ws.Range("L62").Value = 4000
$ws.Range("M62").Value = -3109.3333
$ws.Range("N62").Value = -5248
$ws.Range("H65").Value = 3800
$ws.Range("I65").Value = 3733.3333
$ws.Range("J65").Value = 4000
$ws.Range("K65").Value = 18666.6665
$ws.Range("L65").Value = 20000
$ws.Range("M65").Value = -15546.6665
$ws.Range("N65").Value = -26240

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H31").Value = 0
$ws.Range("J31").Value = 0
$ws.Range("L31").Value = 0
$ws.Range("H68").Value = 553.4211
$ws.Range("I68").Value = 585.625
$ws.Range("K68").Value = 1756.875
$ws.Range("M68").Value = -945.875
$ws.Range("H71").Value = 553.4211
$ws.Range("I71").Value = 585.625
$ws.Range("K71").Value = 5270.625
$ws.Range("M71").Value = -1214.625
$ws.Range("H122").Value = 476.38095
$ws.Range("J122").Value = 640.875
$ws.Range("L122").Value = 5767.875
$ws.Range("N122").Value = -10667.875

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6996.8887
$ws.Range("I70").Value = 5494.25
$ws.Range("J70").Value = 8199
$ws.Range("K70").Value = 5494.25
$ws.Range("L70").Value = 8199
$ws.Range("M70").Value = -5224.25
$ws.Range("N70").Value = -8739
$ws.Range("H73").Value = 6996.8887
$ws.Range("I73").Value = 5494.25
$ws.Range("J73").Value = 8199
$ws.Range("K73").Value = 5494.25
$ws.Range("L73").Value = 8199
$ws.Range("M73").Value = -4558.25
$ws.Range("N73").Value = -10071
$ws.Range("H80").Value = 3590.476
$ws.Range("I80").Value = 3092.2727
$ws.Range("K80").Value = 3092.2727
$ws.Range("M80").Value = -2094.2727
$ws.Range("H83").Value = 3590.476
$ws.Range("I83").Value = 3092.2727
$ws.Range("K83").Value = 15461.3635
$ws.Range("M83").Value = -10469.3635

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 6993.3335
$ws.Range("I68").Value = 0
$ws.Range("K68").Value = 0
$ws.Range("H71").Value = 6993.3335
$ws.Range("I71").Value = 0
$ws.Range("K71").Value = 0
$ws.Range("H132").Value = 49388.08
$ws.Range("I132").Value = 69593.89
$ws.Range("K132").Value = 208781.67
$ws.Range("M132").Value = -206251.67

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 13736.25
$ws.Range("I81").Value = 2122.75
$ws.Range("J81").Value = 25349.75
$ws.Range("K81").Value = 4245.5
$ws.Range("L81").Value = 50699.5
$ws.Range("M81").Value = -3184.5
$ws.Range("N81").Value = -52821.5
$ws.Range("H84").Value = 13736.25
$ws.Range("I84").Value = 2122.75
$ws.Range("J84").Value = 25349.75
$ws.Range("K84").Value = 21227.5
$ws.Range("L84").Value = 253497.5
$ws.Range("M84").Value = -15923.5
$ws.Range("N84").Value = -264105.5

# ---- Delete cells (remove cell entirely, not just clear value) ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("M35").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("N31").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("M68").ClearContents()
$ws.Range("M71").ClearContents()
